# Applies the benchSuite docx update for Renaissance / JDK17 / ShenandoahGC
# (scala-doku, heap-4G): refreshes the summary values at the top of the
# results table and collapses the three tab-separated raw-data rows near
# the bottom down to their single headline values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Top summary rows -----------------------------------------------
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "904"

# Three new rows are inserted right after row 4 (before the row that
# currently holds "0.00004").
$newRow1 = $t.Rows.Add($t.Rows.Item(5))
$newRow1.Cells.Item(1).Range.Text = "0.00002"

$newRow2 = $t.Rows.Add($t.Rows.Item(6))
$newRow2.Cells.Item(1).Range.Text = "0.00073"

$newRow3 = $t.Rows.Add($t.Rows.Item(7))
$newRow3.Cells.Item(1).Range.Text = "0.00015"

# Row 8 ("0.00004") is left untouched; rows 9-12 get new values.
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00026"
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00034"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00036"
$t.Rows.Item(12).Cells.Item(1).Range.Text = "0.17007"

# Rows 13-15 ("0.00018", "0.00022", "0.05866") are removed entirely.
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()

# --- Bottom raw-data rows --------------------------------------------
# These rows hold one run with several <w:t> segments separated by
# <w:tab/> characters; replacing the cell's Range.Text collapses all of
# that down to a single run/value, matching the target edit.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.94"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.17"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "285"
